$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.83636662940318
$ws.Range("C2").Value = 9.858893280965802
$ws.Range("D2").Value = 6.035034581594472
$ws.Range("E2").Value = 12.40552867458618
$ws.Range("F2").Value = 29.17808867001368
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 9.200439912298569
$ws.Range("L2").Value = 9.629963664221867
$ws.Range("O2").Value = 26.16319783465871
$ws.Range("B3").Value = 12.58634283144132
$ws.Range("C3").Value = 9.869088147762509
$ws.Range("D3").Value = 5.999392818380348
$ws.Range("E3").Value = 12.41836552109806
$ws.Range("F3").Value = 29.21694017154243
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 9.019088439655933
$ws.Range("L3").Value = 9.616304712339126
$ws.Range("O3").Value = 26.23728306690385
$ws.Range("B4").Value = 12.43274874063115
$ws.Range("C4").Value = 9.875864392212126
$ws.Range("D4").Value = 5.977044282430189
$ws.Range("E4").Value = 12.42856340197539
$ws.Range("F4").Value = 29.2479740551413
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 8.907323940884739
$ws.Range("L4").Value = 9.609573857461271
$ws.Range("O4").Value = 26.28812220415209
$ws.Range("B5").Value = 12.37022585299062
$ws.Range("C5").Value = 9.878756095433145
$ws.Range("D5").Value = 5.967822980052128
$ws.Range("E5").Value = 12.43330146478733
$ws.Range("F5").Value = 29.26242262540713
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 8.861736647058152
$ws.Range("L5").Value = 9.607249651937085
$ws.Range("O5").Value = 26.31018238476189
$ws.Range("B6").Value = 12.35985067408367
$ws.Range("C6").Value = 9.879244144952214
$ws.Range("D6").Value = 5.966284977222414
$ws.Range("E6").Value = 12.43412338432641
$ws.Range("F6").Value = 29.2649305295869
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 8.854166186231591
$ws.Range("L6").Value = 9.606889068266748
$ws.Range("O6").Value = 26.31392648501471
$ws.Range("B7").Value = 12.43190514288595
$ws.Range("C7").Value = 9.875902862456835
$ws.Range("D7").Value = 5.976920379149378
$ws.Range("E7").Value = 12.4286249434273
$ws.Range("F7").Value = 29.24816162209694
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 8.906709223536961
$ws.Range("L7").Value = 9.609540814374597
$ws.Range("O7").Value = 26.28841428200946
$ws.Range("B8").Value = 12.75022674738786
$ws.Range("C8").Value = 9.862301521174615
$ws.Range("D8").Value = 6.022843271790866
$ws.Range("E8").Value = 12.40947415838358
$ws.Range("F8").Value = 29.18999297856081
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 9.138032400014682
$ws.Range("L8").Value = 9.624911861240317
$ws.Range("O8").Value = 26.18763039704677
$ws.Range("B9").Value = 13.37008440802265
$ws.Range("C9").Value = 9.839708834549315
$ws.Range("D9").Value = 6.1090986815523
$ws.Range("E9").Value = 12.39029198177737
$ws.Range("F9").Value = 29.13299950006157
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 9.585746303022644
$ws.Range("L9").Value = 9.668080173699559
$ws.Range("O9").Value = 26.03255931855462
$ws.Range("B10").Value = 13.8180048508499
$ws.Range("C10").Value = 9.825571311995041
$ws.Range("D10").Value = 6.170004330189494
$ws.Range("E10").Value = 12.38738541562184
$ws.Range("F10").Value = 29.12603353050478
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 9.907714257184464
$ws.Range("L10").Value = 9.707570611229084
$ws.Range("O10").Value = 25.94472382462242
$ws.Range("B11").Value = 14.01919577135675
$ws.Range("C11").Value = 9.819668870438786
$ws.Range("D11").Value = 6.197145035164593
$ws.Range("E11").Value = 12.38848558316909
$ws.Range("F11").Value = 29.13045217163503
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 10.05201181297796
$ws.Range("L11").Value = 9.727183058447169
$ws.Range("O11").Value = 25.91045964479219
$ws.Range("B12").Value = 14.0949378373146
$ws.Range("C12").Value = 9.817509380484198
$ws.Range("D12").Value = 6.207338519019177
$ws.Range("E12").Value = 12.38924960489022
$ws.Range("F12").Value = 29.13321563444181
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 10.10629071537205
$ws.Range("L12").Value = 9.734842849165345
$ws.Range("O12").Value = 25.89830531072873
$ws.Range("B13").Value = 14.07864634749857
$ws.Range("C13").Value = 9.817971107694808
$ws.Range("D13").Value = 6.20514695471318
$ws.Range("E13").Value = 12.38906962471967
$ws.Range("F13").Value = 29.13257200826791
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 10.09461773837186
$ws.Range("L13").Value = 9.733182880733297
$ws.Range("O13").Value = 25.900886427631
$ws.Range("B14").Value = 14.02543644034156
$ws.Range("C14").Value = 9.819489694322774
$ws.Range("D14").Value = 6.197985360758771
$ws.Range("E14").Value = 12.38854148303577
$ws.Range("F14").Value = 29.1306576817252
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 10.05648494053622
$ws.Range("L14").Value = 9.727808588578181
$ws.Range("O14").Value = 25.9094432392505
$ws.Range("B15").Value = 13.99278378964655
$ws.Range("C15").Value = 9.820429710714436
$ws.Range("D15").Value = 6.193587641825768
$ws.Range("E15").Value = 12.3882631935186
$ws.Range("F15").Value = 29.1296270401059
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 10.03307868520305
$ws.Range("L15").Value = 9.72454689739161
$ws.Range("O15").Value = 25.91479148279306
$ws.Range("B16").Value = 13.80479832256583
$ws.Range("C16").Value = 9.825967656502433
$ws.Range("D16").Value = 6.168219001733163
$ws.Range("E16").Value = 12.38736218720226
$ws.Range("F16").Value = 29.12589738847691
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 9.898235962035061
$ws.Range("L16").Value = 9.706321718239069
$ws.Range("O16").Value = 25.94707780981781
$ws.Range("B17").Value = 13.68876229890591
$ws.Range("C17").Value = 9.829500154551701
$ws.Range("D17").Value = 6.152509313564289
$ws.Range("E17").Value = 12.38742928499016
$ws.Range("F17").Value = 29.12555237403577
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 9.814921382407844
$ws.Range("L17").Value = 9.695560499610636
$ws.Range("O17").Value = 25.96834405177636
$ws.Range("B18").Value = 13.62178372926182
$ws.Range("C18").Value = 9.831581749868404
$ws.Range("D18").Value = 6.143420445675778
$ws.Range("E18").Value = 12.3876959514934
$ws.Range("F18").Value = 29.126068207269
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 9.766799947874876
$ws.Range("L18").Value = 9.689526309928615
$ws.Range("O18").Value = 25.98111153046122
$ws.Range("B19").Value = 13.59906746308127
$ws.Range("C19").Value = 9.832295107098227
$ws.Range("D19").Value = 6.140334081574568
$ws.Range("E19").Value = 12.38782543691008
$ws.Range("F19").Value = 29.12636554911616
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 9.750473948707242
$ws.Range("L19").Value = 9.687510041166867
$ws.Range("O19").Value = 25.98552631679534
$ws.Range("B20").Value = 13.70113969140317
$ws.Range("C20").Value = 9.829118963345872
$ws.Range("D20").Value = 6.154187149570886
$ws.Range("E20").Value = 12.38739854360723
$ws.Range("F20").Value = 29.12551517990302
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 9.82381155905215
$ws.Range("L20").Value = 9.696689996785421
$ws.Range("O20").Value = 25.9660247648184
$ws.Range("B21").Value = 14.04107812107484
$ws.Range("C21").Value = 9.819041599053712
$ws.Range("D21").Value = 6.200091199276782
$ws.Range("E21").Value = 12.38868719073213
$ws.Range("F21").Value = 29.13119038868037
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 10.06769572874912
$ws.Range("L21").Value = 9.729380858575448
$ws.Range("O21").Value = 25.90690760450057
$ws.Range("B22").Value = 14.26062459286737
$ws.Range("C22").Value = 9.812896142730164
$ws.Range("D22").Value = 6.229600615043426
$ws.Range("E22").Value = 12.39155381044264
$ws.Range("F22").Value = 29.14125309031472
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 10.2249478748012
$ws.Range("L22").Value = 9.752102189512163
$ws.Range("O22").Value = 25.87305605276004
$ws.Range("B23").Value = 14.14371213556313
$ws.Range("C23").Value = 9.816135897803173
$ws.Range("D23").Value = 6.213896769724679
$ws.Range("E23").Value = 12.38983896413577
$ws.Range("F23").Value = 29.1353015856689
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 10.14123157447197
$ws.Range("L23").Value = 9.739852706195693
$ws.Range("O23").Value = 25.89068477361488
$ws.Range("B24").Value = 13.69554470099804
$ws.Range("C24").Value = 9.82929114188345
$ws.Range("D24").Value = 6.153428777730584
$ws.Range("E24").Value = 12.38741173124351
$ws.Range("F24").Value = 29.12552977065055
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 9.81979300010665
$ws.Range("L24").Value = 9.696178875552166
$ws.Range("O24").Value = 25.96707162868305
$ws.Range("B25").Value = 13.20337514140248
$ws.Range("C25").Value = 9.845386682832444
$ws.Range("D25").Value = 6.086186615065092
$ws.Range("E25").Value = 12.3935150547932
$ws.Range("F25").Value = 29.14229270086278
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 9.465617528423079
$ws.Range("L25").Value = 9.655023561689372
$ws.Range("O25").Value = 26.069937481368
